# Refresh cryptocurrency Price / Volume(1h) columns with the latest scrape values.
# Numeric-looking Price strings are entered with a leading apostrophe so Excel keeps
# them as text (matching the sheet's existing text-formatted Price column), then the
# cell style is restored to Normal so no stray number-format/quote-prefix style sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '41.909.36'
$ws.Cells.Item(2, 5).Value = '  +5.92%  '
$ws.Cells.Item(3, 4).Value = '2.229.76'
$ws.Cells.Item(5, 5).Value = '  +1.98%  '
$ws.Cells.Item(6, 4).Value = '''0.625'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.89%  '
$ws.Cells.Item(7, 4).Value = '''61.06'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -2.83%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 5).Value = '  +2.99%  '
$ws.Cells.Item(10, 4).Value = '''58.96'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.05%  '
$ws.Cells.Item(11, 4).Value = '''0.0893'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +5.00%  '
$ws.Cells.Item(12, 5).Value = '  -0.41%  '
$ws.Cells.Item(13, 4).Value = '2.560.62'
$ws.Cells.Item(13, 5).Value = '  +2.92%  '
$ws.Cells.Item(14, 4).Value = '''15.64'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -1.34%  '
$ws.Cells.Item(15, 4).Value = '''21.74'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.19%  '
$ws.Cells.Item(16, 4).Value = '''0.800'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.85%  '
$ws.Cells.Item(17, 4).Value = '''5.57'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.99%  '
$ws.Cells.Item(18, 4).Value = '2.248.87'
$ws.Cells.Item(18, 5).Value = '  +4.06%  '
$ws.Cells.Item(19, 4).Value = '41.794.62'
$ws.Cells.Item(19, 5).Value = '  +5.62%  '
$ws.Cells.Item(20, 4).Value = '''72.29'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.89%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0892'
$ws.Cells.Item(21, 5).Value = '  -2.09%  '
$ws.Cells.Item(22, 5).Value = '  +0.49%  '
$ws.Cells.Item(23, 4).Value = '''250.67'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +10.11%  '
$ws.Cells.Item(24, 5).Value = '  +0.01%  '
$ws.Cells.Item(25, 5).Value = '  +1.55%  '
$ws.Cells.Item(26, 5).Value = '  -0.23%  '
$ws.Cells.Item(27, 4).Value = '''9.63'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +1.83%  '
$ws.Cells.Item(28, 5).Value = '  +4.12%  '
$ws.Cells.Item(29, 4).Value = '''167.12'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -1.84%  '
$ws.Cells.Item(30, 4).Value = '''19.96'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +1.54%  '
$ws.Cells.Item(31, 5).Value = '  -2.48%  '
$ws.Cells.Item(32, 5).Value = '  -1.43%  '
$ws.Cells.Item(33, 5).Value = '  -0.30%  '
$ws.Cells.Item(34, 5).Value = '  +5.40%  '
$ws.Cells.Item(35, 5).Value = '  +3.30%  '
$ws.Cells.Item(36, 4).Value = '''0.0631'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +2.87%  '
$ws.Cells.Item(37, 4).Value = '''6.63'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -4.62%  '
$ws.Cells.Item(38, 5).Value = '  -2.93%  '
$ws.Cells.Item(39, 4).Value = '''2.35'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -1.51%  '
$ws.Cells.Item(40, 4).Value = '''0.000253'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +28.35%  '
$ws.Cells.Item(41, 5).Value = '  +0.13%  '
$ws.Cells.Item(42, 5).Value = '  +5.96%  '
$ws.Cells.Item(43, 4).Value = '''4.83'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.12%  '
$ws.Cells.Item(44, 4).Value = '''8.59'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +8.75%  '
$ws.Cells.Item(45, 4).Value = '''0.0976'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +6.79%  '
$ws.Cells.Item(46, 4).Value = '''98.91'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.30%  '
$ws.Cells.Item(47, 4).Value = '''1.21'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.78%  '
$ws.Cells.Item(48, 4).Value = '1.472.87'
$ws.Cells.Item(48, 5).Value = '  -2.68%  '
$ws.Cells.Item(49, 4).Value = '''16.50'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -6.85%  '
$ws.Cells.Item(50, 5).Value = '  +0.16%  '
$ws.Cells.Item(51, 5).Value = '  -0.99%  '
